$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column (D) to Text format so numeric-looking values
# (e.g. "0.999", "12.48") are preserved as strings rather than converted to numbers,
# matching the original inlineStr cell type used throughout column D.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "68.281.04"
$ws.Range("E2").Value = "  -0.37%  "

$ws.Range("D3").Value = "3.890.79"
$ws.Range("E3").Value = "  -0.86%  "

$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.09%  "

$ws.Range("D5").Value = "483.59"
$ws.Range("E5").Value = "  -0.09%  "

$ws.Range("D6").Value = "144.97"
$ws.Range("E6").Value = "  -1.96%  "

$ws.Range("D7").Value = "0.623"
$ws.Range("E7").Value = "  +0.46%  "

$ws.Range("E8").Value = "  -0.03%  "

$ws.Range("E9").Value = "  +2.65%  "

$ws.Range("E10").Value = "  +7.50%  "

$ws.Range("D11").Value = "0.0000353"
$ws.Range("E11").Value = "  -0.74%  "

$ws.Range("E12").Value = "  +1.61%  "

$ws.Range("E13").Value = "  +0.25%  "

$ws.Range("D14").Value = "4.504.73"
$ws.Range("E14").Value = "  -0.93%  "

$ws.Range("D15").Value = "3.907.97"
$ws.Range("E15").Value = "  -1.23%  "

$ws.Range("D16").Value = "14.18"
$ws.Range("E16").Value = "  -2.97%  "

$ws.Range("E17").Value = "  -0.62%  "

$ws.Range("D18").Value = "19.93"
$ws.Range("E18").Value = "  +1.14%  "

$ws.Range("E19").Value = "  +0.33%  "

$ws.Range("D20").Value = "68.261.72"
$ws.Range("E20").Value = "  -0.62%  "

$ws.Range("D21").Value = "429.84"
$ws.Range("E21").Value = "  -0.33%  "

$ws.Range("D22").Value = "3.61"
$ws.Range("E22").Value = "  +8.19%  "

$ws.Range("D23").Value = "14.83"
$ws.Range("E23").Value = "  +2.35%  "

$ws.Range("B24").Value = "RenderToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D24").Value = "12.48"
$ws.Range("E24").Value = "  +18.12%  "

$ws.Range("B25").Value = "Litecoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D25").Value = "89.13"
$ws.Range("E25").Value = "  +2.50%  "

$ws.Range("D26").Value = "3.67"
$ws.Range("E26").Value = "  +2.56%  "

$ws.Range("E27").Value = "  -3.03%  "

$ws.Range("D28").Value = "37.31"
$ws.Range("E28").Value = "  -2.13%  "

$ws.Range("D29").Value = "5.68"
$ws.Range("E29").Value = "  -3.50%  "

$ws.Range("D30").Value = "712.78"
$ws.Range("E30").Value = "  +0.77%  "

$ws.Range("D31").Value = "13.50"
$ws.Range("E31").Value = "  +2.10%  "

$ws.Range("E33").Value = "  +3.00%  "

$ws.Range("D34").Value = "61.73"
$ws.Range("E34").Value = "  +5.52%  "

$ws.Range("D35").Value = "0.0₃0877"
$ws.Range("E35").Value = "  -2.24%  "

$ws.Range("E36").Value = "  +10.88%  "

$ws.Range("D37").Value = "41.00"
$ws.Range("E37").Value = "  -0.86%  "

$ws.Range("B38").Value = "TheGraph"
$ws.Range("C38").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D38").Value = "0.399"
$ws.Range("E38").Value = "  +16.24%  "

$ws.Range("B39").Value = "Fetch.AI"
$ws.Range("C39").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D39").Value = "3.05"
$ws.Range("E39").Value = "  +6.36%  "

$ws.Range("E40").Value = "  -2.57%  "

$ws.Range("D41").Value = "0.998"
$ws.Range("E41").Value = "  -0.05%  "

$ws.Range("D42").Value = "0.0498"
$ws.Range("E42").Value = "  +6.50%  "

$ws.Range("D43").Value = "3.08"
$ws.Range("E43").Value = "  +2.83%  "

$ws.Range("E44").Value = "  -3.76%  "

$ws.Range("E45").Value = "  +1.70%  "

$ws.Range("D46").Value = "3.38"
$ws.Range("E46").Value = "  +3.73%  "

$ws.Range("E47").Value = "  -0.10%  "

$ws.Range("D48").Value = "0.0₆0350"
$ws.Range("E48").Value = "  +27.25%  "

$ws.Range("E49").Value = "  -1.14%  "

$ws.Range("D50").Value = "2.12"
$ws.Range("E50").Value = "  -2.07%  "

$ws.Range("D51").Value = "144.63"
$ws.Range("E51").Value = "  -1.65%  "
